# "form now reads data from owl file"
# Appends two new survey response rows (rows 8 and 9) to Sheet1, matching
# the shape of the existing rows 2-7 (A:C and F:R are 0/1 flags, D/E/S/T/U
# are text values stored as shared strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows, keyed by column letter -> value, in column order (A..U) so the
# underlying string table grows the same way it did for the existing rows.
$newRows = @(
    @{
        A = 1; B = 1; C = 1
        D = "fake food"
        E = "real food"
        F = 1; G = 0; H = 0; I = 1; J = 0; K = 1; L = 0; M = 1; N = 1; O = 1; P = 0; Q = 1; R = 1
        S = "Moderate"
        T = "amsterdam"
        U = "amsterdamNoord"
    },
    @{
        A = 1; B = 1; C = 0
        D = "fake food 2"
        E = "real food 2"
        F = 1; G = 0; H = 0; I = 1; J = 0; K = 1; L = 0; M = 1; N = 1; O = 0; P = 0; Q = 1; R = 0
        S = "Cheap"
        T = "innsbruck"
        U = "reichenau"
    }
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U")

$startRow = 8
for ($col = 0; $col -lt $columns.Count; $col++) {
    $colLetter = $columns[$col]
    for ($i = 0; $i -lt $newRows.Count; $i++) {
        $rowNum = $startRow + $i
        $ws.Range("$colLetter$rowNum").Value = $newRows[$i][$colLetter]
    }
}
